# Generate Report for Handback
# Updates the handback-status workbook with a fresh run's file identifiers
# (GUIDs), content hashes, and timestamps across the Overview / zh-cn / de-de
# sheets, mirroring a new CI "handback status" report generation.

$wb = $excel.ActiveWorkbook

$missing = [System.Reflection.Missing]::Value

# ---------------------------------------------------------------------------
# New identifiers / timestamps for this report run
# ---------------------------------------------------------------------------
$file1Guid = "dc6764d4-22ed-4c83-b542-7d77762fb365"
$file2Guid = "ffffe48b1004-579a-4ca6-90fc-1cd553717b74"
$xlfHash   = "e500923cea549843464982ed83e6c7f083fa8a77"

$file1Name = "$file1Guid.md"
$file2Name = "$file2Guid.md"
$file1Path = "e2e\$file1Guid.md"
$file2Path = "e2e\$file2Guid.md"

$latestHoDate = "2016-09-02 21:12:30"

$zhHandoffDate  = "2016-09-02 21:12:24"
$zhHandbackDate = "2016-09-02 21:12:42"
$deHandoffDate  = $latestHoDate
$deHandbackDate = "2016-09-02 21:12:50"

# Both rows now point at the same (newly regenerated) xlf pair, since the
# second source file collapsed onto the first file's handoff/handback pass.
$zhXlfName = "$file1Guid.$xlfHash.zh-cn.xlf"
$deXlfName = "$file1Guid.$xlfHash.de-de.xlf"

# ---------------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$ov_rId2 = $wsOverview.Hyperlinks.Item(1).Address
$ov_rId3 = $wsOverview.Hyperlinks.Item(2).Address
if ([string]::IsNullOrEmpty($ov_rId2)) {
    $ov_rId2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0cecea3488614d24e9b08e2b39ab45498c33d1ba/e2e/28f863b0-2b82-4bf4-81d6-4c8c79647f17.md"
}
if ([string]::IsNullOrEmpty($ov_rId3)) {
    $ov_rId3 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0cecea3488614d24e9b08e2b39ab45498c33d1ba/e2e/b66ee5bc-4b79-4695-8358-74fdbabe6e56.md"
}

$wsOverview.Range("A2").Value = $file1Name
$wsOverview.Range("A3").Value = $file2Name
$wsOverview.Range("G2").Value = $latestHoDate
$wsOverview.Range("G3").Value = $latestHoDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $ov_rId2, $missing, $missing, $file1Path) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $ov_rId3, $missing, $missing, $file2Path) | Out-Null

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zh_rId2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0cecea3488614d24e9b08e2b39ab45498c33d1ba/e2e/28f863b0-2b82-4bf4-81d6-4c8c79647f17.md"
$zh_rId3 = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/76642b8b5ff2036a299c4ef63aef55fd23f65421/e2e/28f863b0-2b82-4bf4-81d6-4c8c79647f17.md"
$zh_rId4 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0cecea3488614d24e9b08e2b39ab45498c33d1ba/e2e/b66ee5bc-4b79-4695-8358-74fdbabe6e56.md"
$zh_rId5 = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/76642b8b5ff2036a299c4ef63aef55fd23f65421/e2e/b66ee5bc-4b79-4695-8358-74fdbabe6e56.md"

$wsZh.Range("A2").Value = $file1Name
$wsZh.Range("G2").Value = $zhXlfName
$wsZh.Range("H2").Value = $zhHandoffDate
$wsZh.Range("I2").Value = $file1Name
$wsZh.Range("J2").Value = $zhXlfName
$wsZh.Range("K2").Value = $zhHandbackDate

$wsZh.Range("A3").Value = $file2Name
$wsZh.Range("G3").Value = $zhXlfName
$wsZh.Range("H3").Value = $zhHandoffDate
$wsZh.Range("I3").Value = $file2Name
$wsZh.Range("J3").Value = $zhXlfName
$wsZh.Range("K3").Value = $zhHandbackDate

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $zh_rId2, $missing, $missing, $file1Name) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $zh_rId3, $missing, $missing, $file1Name) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $zh_rId4, $missing, $missing, $file2Name) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $zh_rId5, $missing, $missing, $file2Name) | Out-Null

# ---------------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$de_rId2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0cecea3488614d24e9b08e2b39ab45498c33d1ba/e2e/28f863b0-2b82-4bf4-81d6-4c8c79647f17.md"
$de_rId3 = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/0d51b400ab38340b61772736579b5d6de68f15c4/e2e/28f863b0-2b82-4bf4-81d6-4c8c79647f17.md"
$de_rId4 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0cecea3488614d24e9b08e2b39ab45498c33d1ba/e2e/b66ee5bc-4b79-4695-8358-74fdbabe6e56.md"
$de_rId5 = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/0d51b400ab38340b61772736579b5d6de68f15c4/e2e/b66ee5bc-4b79-4695-8358-74fdbabe6e56.md"

$wsDe.Range("A2").Value = $file1Name
$wsDe.Range("G2").Value = $deXlfName
$wsDe.Range("H2").Value = $deHandoffDate
$wsDe.Range("I2").Value = $file1Name
$wsDe.Range("J2").Value = $deXlfName
$wsDe.Range("K2").Value = $deHandbackDate

$wsDe.Range("A3").Value = $file2Name
$wsDe.Range("G3").Value = $deXlfName
$wsDe.Range("H3").Value = $deHandoffDate
$wsDe.Range("I3").Value = $file2Name
$wsDe.Range("J3").Value = $deXlfName
$wsDe.Range("K3").Value = $deHandbackDate

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $de_rId2, $missing, $missing, $file1Name) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $de_rId3, $missing, $missing, $file1Name) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $de_rId4, $missing, $missing, $file2Name) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $de_rId5, $missing, $missing, $file2Name) | Out-Null
